# Auto-generated edit script: updates market-price derived columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# on several rows across multiple sheets, reflecting refreshed market data
# pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3179.4
$ws.Cells.Item(76, 9).Value = 2851.8948
$ws.Cells.Item(76, 10).Value = 4216.5
$ws.Cells.Item(76, 11).Value = 2851.8948
$ws.Cells.Item(76, 12).Value = 4216.5
$ws.Cells.Item(76, 13).Value = -2536.8948
$ws.Cells.Item(76, 14).Value = -4846.5
$ws.Cells.Item(79, 8).Value = 3179.4
$ws.Cells.Item(79, 9).Value = 2851.8948
$ws.Cells.Item(79, 10).Value = 4216.5
$ws.Cells.Item(79, 11).Value = 2851.8948
$ws.Cells.Item(79, 12).Value = 4216.5
$ws.Cells.Item(79, 13).Value = -1759.8948
$ws.Cells.Item(79, 14).Value = -6400.5
$ws.Cells.Item(80, 8).Value = 3567.75
$ws.Cells.Item(80, 9).Value = 403.3
$ws.Cells.Item(80, 10).Value = 5828.0713
$ws.Cells.Item(80, 11).Value = 1209.9
$ws.Cells.Item(80, 12).Value = 17484.2139
$ws.Cells.Item(80, 13).Value = -211.9000000000001
$ws.Cells.Item(80, 14).Value = -19480.2139
$ws.Cells.Item(83, 8).Value = 3567.75
$ws.Cells.Item(83, 9).Value = 403.3
$ws.Cells.Item(83, 10).Value = 5828.0713
$ws.Cells.Item(83, 11).Value = 3629.7
$ws.Cells.Item(83, 12).Value = 52452.64169999999
$ws.Cells.Item(83, 13).Value = 1362.3
$ws.Cells.Item(83, 14).Value = -62436.64169999999
$ws.Cells.Item(129, 8).Value = 3379541.5
$ws.Cells.Item(129, 10).Value = 1161.2239
$ws.Cells.Item(129, 12).Value = 3483.6717
$ws.Cells.Item(129, 14).Value = -13483.6717
$ws.Cells.Item(132, 8).Value = 3451069.5
$ws.Cells.Item(132, 9).Value = 4350209.5
$ws.Cells.Item(132, 10).Value = 4366.4165
$ws.Cells.Item(132, 11).Value = 13050628.5
$ws.Cells.Item(132, 12).Value = 13099.2495
$ws.Cells.Item(132, 13).Value = -13048098.5
$ws.Cells.Item(132, 14).Value = -18159.2495
$ws.Cells.Item(137, 8).Value = 1726395.1
$ws.Cells.Item(137, 9).Value = 2502381.5
$ws.Cells.Item(137, 10).Value = 1980.7222
$ws.Cells.Item(137, 11).Value = 7507144.5
$ws.Cells.Item(137, 12).Value = 5942.1666
$ws.Cells.Item(137, 13).Value = -7504594.5
$ws.Cells.Item(137, 14).Value = -11042.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1321.6316
$ws.Cells.Item(45, 9).Value = 1040.4
$ws.Cells.Item(45, 10).Value = 2376.25
$ws.Cells.Item(45, 11).Value = 1040.4
$ws.Cells.Item(45, 12).Value = 2376.25
$ws.Cells.Item(45, 13).Value = -663.4000000000001
$ws.Cells.Item(45, 14).Value = -3130.25
$ws.Cells.Item(61, 8).Value = 1483.091
$ws.Cells.Item(61, 9).Value = 691.5
$ws.Cells.Item(61, 10).Value = 4040.5386
$ws.Cells.Item(61, 11).Value = 691.5
$ws.Cells.Item(61, 12).Value = 4040.5386
$ws.Cells.Item(61, 13).Value = -479.5
$ws.Cells.Item(61, 14).Value = -4464.5386
$ws.Cells.Item(122, 8).Value = 3022.318
$ws.Cells.Item(122, 9).Value = 2566.0667
$ws.Cells.Item(122, 10).Value = 4000
$ws.Cells.Item(122, 11).Value = 7698.2001
$ws.Cells.Item(122, 12).Value = 12000
$ws.Cells.Item(122, 13).Value = -5248.2001
$ws.Cells.Item(122, 14).Value = -16900
$ws.Cells.Item(123, 8).Value = 29899.5
$ws.Cells.Item(123, 10).Value = 29899.5
$ws.Cells.Item(123, 12).Value = 29899.5
$ws.Cells.Item(123, 14).Value = -39699.5
$ws.Cells.Item(132, 8).Value = 1718.6482
$ws.Cells.Item(132, 9).Value = 1250.8298
$ws.Cells.Item(132, 10).Value = 4859.7144
$ws.Cells.Item(132, 11).Value = 3752.4894
$ws.Cells.Item(132, 12).Value = 14579.1432
$ws.Cells.Item(132, 13).Value = -1222.4894
$ws.Cells.Item(132, 14).Value = -19639.1432
$ws.Cells.Item(134, 8).Value = 36233.332
$ws.Cells.Item(134, 10).Value = 36233.332
$ws.Cells.Item(134, 12).Value = 36233.332
$ws.Cells.Item(134, 14).Value = -46373.332
$ws.Cells.Item(136, 8).Value = 1483.091
$ws.Cells.Item(136, 9).Value = 691.5
$ws.Cells.Item(136, 10).Value = 4040.5386
$ws.Cells.Item(136, 11).Value = 2074.5
$ws.Cells.Item(136, 12).Value = 12121.6158
$ws.Cells.Item(136, 13).Value = 475.5
$ws.Cells.Item(136, 14).Value = -17221.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 3567.1333
$ws.Cells.Item(99, 9).Value = 2504.5
$ws.Cells.Item(99, 10).Value = 4781.5713
$ws.Cells.Item(99, 11).Value = 2504.5
$ws.Cells.Item(99, 12).Value = 4781.5713
$ws.Cells.Item(99, 13).Value = -1006.5
$ws.Cells.Item(99, 14).Value = -7777.5713
$ws.Cells.Item(134, 8).Value = 1571.8966
$ws.Cells.Item(134, 9).Value = 704.6
$ws.Cells.Item(134, 10).Value = 6992.5
$ws.Cells.Item(134, 11).Value = 2113.8
$ws.Cells.Item(134, 12).Value = 20977.5
$ws.Cells.Item(134, 13).Value = 421.1999999999998
$ws.Cells.Item(134, 14).Value = -26047.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 9261762
$ws.Cells.Item(58, 9).Value = 1482
$ws.Cells.Item(58, 10).Value = 41672744
$ws.Cells.Item(58, 11).Value = 1482
$ws.Cells.Item(58, 12).Value = 41672744
$ws.Cells.Item(58, 13).Value = -1279
$ws.Cells.Item(58, 14).Value = -41673150
$ws.Cells.Item(132, 8).Value = 2189.6875
$ws.Cells.Item(132, 9).Value = 1573.9286
$ws.Cells.Item(132, 10).Value = 6500
$ws.Cells.Item(132, 11).Value = 4721.7858
$ws.Cells.Item(132, 12).Value = 19500
$ws.Cells.Item(132, 13).Value = -2191.7858
$ws.Cells.Item(132, 14).Value = -24560
$ws.Cells.Item(134, 8).Value = 2090.3704
$ws.Cells.Item(134, 9).Value = 1042.091
$ws.Cells.Item(134, 10).Value = 6702.8
$ws.Cells.Item(134, 11).Value = 3126.273
$ws.Cells.Item(134, 12).Value = 20108.4
$ws.Cells.Item(134, 13).Value = -591.2729999999997
$ws.Cells.Item(134, 14).Value = -25178.4
$ws.Cells.Item(136, 8).Value = 9261762
$ws.Cells.Item(136, 9).Value = 1482
$ws.Cells.Item(136, 10).Value = 41672744
$ws.Cells.Item(136, 11).Value = 4446
$ws.Cells.Item(136, 12).Value = 125018232
$ws.Cells.Item(136, 13).Value = -1896
$ws.Cells.Item(136, 14).Value = -125023332

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1724.3334
$ws.Cells.Item(5, 9).Value = 1182.4
$ws.Cells.Item(5, 10).Value = 2111.4285
$ws.Cells.Item(5, 11).Value = 3547.2
$ws.Cells.Item(5, 12).Value = 6334.2855
$ws.Cells.Item(5, 13).Value = -3435.2
$ws.Cells.Item(5, 14).Value = -6558.2855
$ws.Cells.Item(18, 8).Value = 740
$ws.Cells.Item(18, 9).Value = 372.5
$ws.Cells.Item(18, 11).Value = 1117.5
$ws.Cells.Item(18, 13).Value = -948.5
$ws.Cells.Item(57, 8).Value = 3556
$ws.Cells.Item(57, 10).Value = 4000
$ws.Cells.Item(57, 12).Value = 12000
$ws.Cells.Item(57, 14).Value = -13118
$ws.Cells.Item(64, 8).Value = 2288.7778
$ws.Cells.Item(64, 9).Value = 867
$ws.Cells.Item(64, 10).Value = 2999.6667
$ws.Cells.Item(64, 11).Value = 2601
$ws.Cells.Item(64, 12).Value = 8999.000100000001
$ws.Cells.Item(64, 13).Value = -2331
$ws.Cells.Item(64, 14).Value = -9539.000100000001
$ws.Cells.Item(67, 8).Value = 2288.7778
$ws.Cells.Item(67, 9).Value = 867
$ws.Cells.Item(67, 10).Value = 2999.6667
$ws.Cells.Item(67, 11).Value = 2601
$ws.Cells.Item(67, 12).Value = 8999.000100000001
$ws.Cells.Item(67, 13).Value = -1665
$ws.Cells.Item(67, 14).Value = -10871.0001
$ws.Cells.Item(116, 8).Value = 1850
$ws.Cells.Item(116, 9).Value = 1133.3334
$ws.Cells.Item(116, 10).Value = 4000
$ws.Cells.Item(116, 11).Value = 3400.0002
$ws.Cells.Item(116, 12).Value = 12000
$ws.Cells.Item(116, 13).Value = 41.99980000000005
$ws.Cells.Item(116, 14).Value = -18884
$ws.Cells.Item(118, 8).Value = 2371.9
$ws.Cells.Item(118, 9).Value = 999.6667
$ws.Cells.Item(118, 11).Value = 2999.0001
$ws.Cells.Item(118, 13).Value = -1756.0001
$ws.Cells.Item(122, 8).Value = 928.1579
$ws.Cells.Item(122, 9).Value = 451.22223
$ws.Cells.Item(122, 10).Value = 1357.4
$ws.Cells.Item(122, 11).Value = 4061.00007
$ws.Cells.Item(122, 12).Value = 12216.6
$ws.Cells.Item(122, 13).Value = -1611.00007
$ws.Cells.Item(122, 14).Value = -17116.6
$ws.Cells.Item(123, 8).Value = 2802.8667
$ws.Cells.Item(123, 9).Value = 455
$ws.Cells.Item(123, 10).Value = 3164.077
$ws.Cells.Item(123, 11).Value = 1365
$ws.Cells.Item(123, 12).Value = 9492.231
$ws.Cells.Item(123, 13).Value = 1085
$ws.Cells.Item(123, 14).Value = -14392.231
$ws.Cells.Item(131, 8).Value = 1563.0488
$ws.Cells.Item(131, 9).Value = 2616
$ws.Cells.Item(131, 10).Value = 1223.3871
$ws.Cells.Item(131, 11).Value = 7848
$ws.Cells.Item(131, 12).Value = 3670.1613
$ws.Cells.Item(131, 13).Value = -2808
$ws.Cells.Item(131, 14).Value = -13750.1613
$ws.Cells.Item(135, 8).Value = 1724.3334
$ws.Cells.Item(135, 9).Value = 1182.4
$ws.Cells.Item(135, 10).Value = 2111.4285
$ws.Cells.Item(135, 11).Value = 10641.6
$ws.Cells.Item(135, 12).Value = 19002.8565
$ws.Cells.Item(135, 13).Value = -8106.6
$ws.Cells.Item(135, 14).Value = -24072.8565

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2358.75
$ws.Cells.Item(40, 9).Value = 1466.6666
$ws.Cells.Item(40, 10).Value = 3250.8333
$ws.Cells.Item(40, 11).Value = 1466.6666
$ws.Cells.Item(40, 12).Value = 3250.8333
$ws.Cells.Item(40, 13).Value = -1330.6666
$ws.Cells.Item(40, 14).Value = -3522.8333
$ws.Cells.Item(46, 8).Value = 1953.4
$ws.Cells.Item(46, 9).Value = 520.2
$ws.Cells.Item(46, 10).Value = 2670
$ws.Cells.Item(46, 11).Value = 520.2
$ws.Cells.Item(46, 12).Value = 2670
$ws.Cells.Item(46, 13).Value = -332.2
$ws.Cells.Item(46, 14).Value = -3046
$ws.Cells.Item(122, 8).Value = 3053.6
$ws.Cells.Item(122, 9).Value = 2650.3333
$ws.Cells.Item(122, 10).Value = 4666.6665
$ws.Cells.Item(122, 11).Value = 7950.999899999999
$ws.Cells.Item(122, 12).Value = 13999.9995
$ws.Cells.Item(122, 13).Value = -5500.999899999999
$ws.Cells.Item(122, 14).Value = -18899.9995
$ws.Cells.Item(127, 8).Value = 30000
$ws.Cells.Item(127, 10).Value = 30000
$ws.Cells.Item(127, 12).Value = 30000
$ws.Cells.Item(127, 14).Value = -39920
$ws.Cells.Item(132, 8).Value = 1990.079
$ws.Cells.Item(132, 9).Value = 1249.3704
$ws.Cells.Item(132, 10).Value = 3808.182
$ws.Cells.Item(132, 11).Value = 3748.1112
$ws.Cells.Item(132, 12).Value = 11424.546
$ws.Cells.Item(132, 13).Value = -1218.1112
$ws.Cells.Item(132, 14).Value = -16484.546
$ws.Cells.Item(136, 8).Value = 2944379.5
$ws.Cells.Item(136, 9).Value = 4168585.5
$ws.Cells.Item(136, 10).Value = 6285.5
$ws.Cells.Item(136, 11).Value = 12505756.5
$ws.Cells.Item(136, 12).Value = 18856.5
$ws.Cells.Item(136, 13).Value = -12503206.5
$ws.Cells.Item(136, 14).Value = -23956.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 591036.75
$ws.Cells.Item(122, 9).Value = 1252153.2
$ws.Cells.Item(122, 10).Value = 3377.6667
$ws.Cells.Item(122, 11).Value = 3756459.6
$ws.Cells.Item(122, 12).Value = 10133.0001
$ws.Cells.Item(122, 13).Value = -3754009.6
$ws.Cells.Item(122, 14).Value = -15033.0001
$ws.Cells.Item(132, 8).Value = 286871.75
$ws.Cells.Item(132, 9).Value = 438423.78
$ws.Cells.Item(132, 10).Value = 37893.43
$ws.Cells.Item(132, 11).Value = 1315271.34
$ws.Cells.Item(132, 12).Value = 113680.29
$ws.Cells.Item(132, 13).Value = -1312741.34
$ws.Cells.Item(132, 14).Value = -118740.29
$ws.Cells.Item(136, 8).Value = 1034.5111
$ws.Cells.Item(136, 9).Value = 602.1818
$ws.Cells.Item(136, 10).Value = 1448.0435
$ws.Cells.Item(136, 11).Value = 1806.5454
$ws.Cells.Item(136, 12).Value = 4344.1305
$ws.Cells.Item(136, 13).Value = 743.4546
$ws.Cells.Item(136, 14).Value = -9444.130499999999

